$wb = $excel.ActiveWorkbook

# --- Add new worksheet "9_" after the last sheet ("8_") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "9_"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 18.5703125
$ws.Columns.Item(2).ColumnWidth = 21
$ws.Columns.Item(3).ColumnWidth = 33.5703125
$ws.Columns.Item(4).ColumnWidth = 31

# --- Cell values ---
$ws.Range("A1").Value = 'Look at the reparameterized equation, and match the condition below with the behavior of the model.  What happens when:'
$ws.Range("B1").Value = 'Correct order of definitions'
$ws.Range("C1").Value = 'Definitions'
$ws.Range("A2").Value = 'p'' is much smaller than ''K'''
$ws.Range("B2").Value = 'B'
$ws.Range("C2").Value = 'Delta p will be large and negative: the population will shrink'
$ws.Range("D2").Value = 'The rate of growth will depend on ''r'''
$ws.Range("A3").Value = 'p'' is close to the value of ''K'''
$ws.Range("B3").Value = 'C'
$ws.Range("C3").Value = 'Delta p will tend to be large; its sign will cause the population to move quickly to equilibrium'
$ws.Range("D3").Value = ' '
$ws.Range("A4").Value = 'p'' is much larger than K'
$ws.Range("B4").Value = 'A'
$ws.Range("C4").Value = 'Delta p will be small not matter what the other variables are'
$ws.Range("A5").Value = 'r'' is very large'
$ws.Range("B5").Value = 'E'
$ws.Range("C5").Value = 'Delta p will tend to be large and positive, which will cause the population to grow quickly'
$ws.Range("A6").Value = 'r'' is very small'
$ws.Range("B6").Value = 'C'
$ws.Range("C6").Value = 'Delta p will be positive but the rate of growth will depend on other variables'

# --- Alignment / wrap text styling ---
# style 1: wrap text only
$ws.Range("C1:C6").WrapText = $true
$ws.Range("D1:D3").WrapText = $true
# style 2: center horizontal+vertical, wrap text
$ws.Range("B1:B6").WrapText = $true
$ws.Range("B1:B6").HorizontalAlignment = -4108
$ws.Range("B1:B6").VerticalAlignment = -4108
# style 3: wrap text (quote-prefixed originally, text already correct)
$ws.Range("A1:A6").WrapText = $true

# --- D4:D6 are present but empty in the source; touch them so the dimension/style matches ---
$ws.Range("D1,D4,D5,D6").WrapText = $true

# --- Sheet tab/selection state: new sheet is active, D3 selected ---
$ws.Range("D3").Select()

# --- Previously-active sheet "8_" now shows the whole data range selected, no single active cell ---
$prev = $wb.Worksheets.Item("8_")
$prev.Range("A1:D6").Select()

# --- Reactivate the new sheet so it ends up as the active/selected tab ---
$ws.Activate()
$ws.Range("D3").Select()
